$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.402.72"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "1.803.42"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.68"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.580"
$ws.Range("E6").Value = "  +4.08%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "35.05"
$ws.Range("E8").Value = "  +6.61%  "
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").Value = "2.062.63"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.812.03"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.16"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").Value = "34.398.16"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.93"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.09"
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0795"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "170.40"
$ws.Range("E24").Value = "  +3.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.10"
$ws.Range("E25").Value = "  +2.49%  "
$ws.Range("E26").Value = "  +4.04%  "
$ws.Range("E27").Value = "  +0.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.118"
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.96"
$ws.Range("E30").Value = "  -4.24%  "
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("D35").Value = "1.396.90"
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.679"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.52"
$ws.Range("E37").Value = "  -2.30%  "
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0190"
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.82"
$ws.Range("E40").Value = "  -2.99%  "
$ws.Range("E41").Value = "  +3.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.947"
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.39"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("E45").Value = "  +2.67%  "
$ws.Range("E46").Value = "  -2.49%  "
$ws.Range("E47").Value = "  -1.37%  "
$ws.Range("D48").Value = "1.962.94"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.47"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D51").Value = "0.0₆0130"
$ws.Range("E51").Value = "  +0.75%  "
